$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'97.699.70"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.51%  "
$ws.Range("D3").Value = "'3.403.33"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.29%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'255.06"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.07%  "
$ws.Range("D6").Value = "'649.89"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.95%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").Value = "'0.425"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.60%  "
$ws.Range("E9").Value = "  +6.87%  "
$ws.Range("E10").Value = "  -0.02%  "
$ws.Range("D11").Value = "'3.400.08"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.22%  "
$ws.Range("E12").Value = "  +4.83%  "
$ws.Range("D13").Value = "'41.23"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.87%  "
$ws.Range("D14").Value = "'6.23"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +13.08%  "
$ws.Range("E15").Value = "  +2.28%  "
$ws.Range("D16").Value = "'97.311.86"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.56%  "
$ws.Range("D17").Value = "'4.029.68"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.01%  "
$ws.Range("D18").Value = "'8.46"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +31.24%  "
$ws.Range("D19").Value = "'3.401.50"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.32%  "
$ws.Range("D20").Value = "'17.28"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +12.83%  "
$ws.Range("B21").Value = "Stellar"
$ws.Range("C21").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D21").Value = "'0.486"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +39.27%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "'10.68"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +14.15%  "
$ws.Range("D23").Value = "'3.41"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.89%  "
$ws.Range("D24").Value = "'499.15"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.67%  "
$ws.Range("E25").Value = "  +0.16%  "
$ws.Range("D26").Value = "'6.12"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +7.74%  "
$ws.Range("D27").Value = "'98.16"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +9.98%  "
$ws.Range("D28").Value = "'12.54"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.77%  "
$ws.Range("D29").Value = "'3.583.66"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.52%  "
$ws.Range("E30").Value = "  +9.71%  "
$ws.Range("E31").Value = "  +5.22%  "
$ws.Range("E32").Value = "  -0.13%  "
$ws.Range("D33").Value = "'11.20"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.55%  "
$ws.Range("D34").Value = "'1.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.06%  "
$ws.Range("D35").Value = "'0.565"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Value = "'29.43"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.97%  "
$ws.Range("E37").Value = "  +14.34%  "
$ws.Range("D38").Value = "'7.65"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.49%  "
$ws.Range("E39").Value = "  +0.72%  "
$ws.Range("E40").Value = "  +12.83%  "
$ws.Range("D41").Value = "'509.62"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.63%  "
$ws.Range("E42").Value = "  -0.11%  "
$ws.Range("E43").Value = "  +10.15%  "
$ws.Range("E44").Value = "  -3.79%  "
$ws.Range("D45").Value = "'0.0413"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +21.79%  "
$ws.Range("D46").Value = "'5.46"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +14.00%  "
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("E48").Value = "  +2.73%  "
$ws.Range("D49").Value = "'8.14"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +10.49%  "
$ws.Range("D50").Value = "'1.55"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +12.54%  "
$ws.Range("E51").Value = "  +10.83%  "
